$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, centered, bordered) from the existing
# header cell H1 onto the two new header cells before setting their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for new columns I (I0) and J (IF), rows 2-26
$data = @(
    @(1, 4),   # row 2
    @(6, 8),   # row 3
    @(3, 7),   # row 4
    @(1, 7),   # row 5
    @(1, 5),   # row 6
    @(1, 6),   # row 7
    @(1, 6),   # row 8
    @(5, 8),   # row 9
    @(1, 6),   # row 10
    @(1, 6),   # row 11
    @(1, 6),   # row 12
    @(1, 6),   # row 13
    @(1, 7),   # row 14
    @(1, 6),   # row 15
    @(1, 5),   # row 16
    @(1, 1),   # row 17
    @(1, 4),   # row 18
    @(1, 4),   # row 19
    @(1, 5),   # row 20
    @(1, 5),   # row 21
    @(1, 5),   # row 22
    @(1, 5),   # row 23
    @(1, 3),   # row 24
    @(7, 9),   # row 25
    @(1, 3)    # row 26
)

$r = 2
foreach ($pair in $data) {
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
    $r++
}
